$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H55").Value = 2415.0833
$ws.Range("I55").Value = 295
$ws.Range("J55").Value = 3121.7778
$ws.Range("K55").Value = 295
$ws.Range("L55").Value = 3121.7778
$ws.Range("M55").Value = -81
$ws.Range("N55").Value = -3549.7778

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H58").Value = 381.5
$ws.Range("I58").Value = 381.5
$ws.Range("K58").Value = 1144.5
$ws.Range("M58").Value = -994.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H70").Value = 1709.7
$ws.Range("I70").Value = 1683.9231
$ws.Range("J70").Value = 1757.5714
$ws.Range("K70").Value = 5051.7693
$ws.Range("L70").Value = 5272.7142
$ws.Range("M70").Value = -4781.7693
$ws.Range("N70").Value = -5812.7142

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H73").Value = 1709.7
$ws.Range("I73").Value = 1683.9231
$ws.Range("J73").Value = 1757.5714
$ws.Range("K73").Value = 5051.7693
$ws.Range("L73").Value = 5272.7142
$ws.Range("M73").Value = -4115.7693
$ws.Range("N73").Value = -7144.7142

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H118").Value = 682.625
$ws.Range("I118").Value = 709.7143
$ws.Range("K118").Value = 2129.1429
$ws.Range("M118").Value = -472.1428999999998

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H121").Value = 2800
$ws.Range("J121").Value = 2800
$ws.Range("L121").Value = 8400
$ws.Range("N121").Value = -11894

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 2360.6
$ws.Range("I132").Value = 2031.4286
$ws.Range("K132").Value = 6094.2858
$ws.Range("M132").Value = -3564.2858

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 6350.185
$ws.Range("I137").Value = 3676.7727
$ws.Range("K137").Value = 11030.3181
$ws.Range("M137").Value = -8480.3181

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H37").Value = 49000
$ws.Range("I37").Value = 0
$ws.Range("K37").Value = 0
$ws.Range("M37").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 7226129.5
$ws.Range("I74").Value = 11366163
$ws.Range("K74").Value = 11366163
$ws.Range("M74").Value = -11365289

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 7226129.5
$ws.Range("I77").Value = 11366163
$ws.Range("K77").Value = 56830815
$ws.Range("M77").Value = -56826447

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H102").Value = 5560.476
$ws.Range("I102").Value = 6547
$ws.Range("J102").Value = 1943.2222
$ws.Range("K102").Value = 6547
$ws.Range("L102").Value = 1943.2222
$ws.Range("M102").Value = -4925
$ws.Range("N102").Value = -5187.2222

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 3781.6667
$ws.Range("I132").Value = 2362.3462
$ws.Range("J132").Value = 5723.8945
$ws.Range("K132").Value = 7087.0386
$ws.Range("L132").Value = 17171.6835
$ws.Range("M132").Value = -4557.0386
$ws.Range("N132").Value = -22231.6835

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H139").Value = 67499.5
$ws.Range("J139").Value = 67499.5
$ws.Range("L139").Value = 67499.5
$ws.Range("N139").Value = -77779.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 295262.38
$ws.Range("I134").Value = 1203.6
$ws.Range("J134").Value = 2500703.2
$ws.Range("K134").Value = 3610.8
$ws.Range("L134").Value = 7502109.600000001
$ws.Range("M134").Value = -1075.8
$ws.Range("N134").Value = -7507179.600000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2088.4285
$ws.Range("I16").Value = 2266.6667
$ws.Range("K16").Value = 2266.6667
$ws.Range("M16").Value = -1979.6667

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H86").Value = 5848.125
$ws.Range("I86").Value = 5541
$ws.Range("K86").Value = 5541
$ws.Range("M86").Value = -4418

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H89").Value = 5848.125
$ws.Range("I89").Value = 5541
$ws.Range("K89").Value = 27705
$ws.Range("M89").Value = -22089

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H113").Value = 2088.4285
$ws.Range("I113").Value = 2266.6667
$ws.Range("K113").Value = 2266.6667
$ws.Range("M113").Value = -96.66670000000022

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 385.5
$ws.Range("I122").Value = 385.5
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 1156.5
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = 1293.5
$ws.Range("N122").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 1810.3
$ws.Range("I132").Value = 1638
$ws.Range("J132").Value = 2499.5
$ws.Range("K132").Value = 4914
$ws.Range("L132").Value = 7498.5
$ws.Range("M132").Value = -2384
$ws.Range("N132").Value = -12558.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 2826.25
$ws.Range("I134").Value = 1651.7667
$ws.Range("K134").Value = 4955.300099999999
$ws.Range("M134").Value = -2420.300099999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H14").Value = 419.16666
$ws.Range("I14").Value = 419.16666
$ws.Range("K14").Value = 1257.49998
$ws.Range("M14").Value = -1084.49998

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 5248.75
$ws.Range("J34").Value = 5248.75
$ws.Range("L34").Value = 15746.25
$ws.Range("N34").Value = -15914.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H52").Value = 0
$ws.Range("J52").Value = 0
$ws.Range("L52").Value = 0
$ws.Range("N52").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1496.619
$ws.Range("I102").Value = 1120.2667
$ws.Range("K102").Value = 1120.2667
$ws.Range("M102").Value = 501.7333000000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 2209.818
$ws.Range("I122").Value = 1917.1666
$ws.Range("K122").Value = 5751.4998
$ws.Range("M122").Value = -3301.4998

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 4999.5
$ws.Range("I126").Value = 4999.5
$ws.Range("K126").Value = 14998.5
$ws.Range("M126").Value = -12528.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1624.75
$ws.Range("I22").Value = 1500
$ws.Range("J22").Value = 1666.3334
$ws.Range("K22").Value = 1500
$ws.Range("L22").Value = 1666.3334
$ws.Range("M22").Value = -1205
$ws.Range("N22").Value = -2256.3334

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H27").Value = 1624.75
$ws.Range("I27").Value = 1500
$ws.Range("J27").Value = 1666.3334
$ws.Range("K27").Value = 1500
$ws.Range("L27").Value = 1666.3334
$ws.Range("M27").Value = -1393
$ws.Range("N27").Value = -1880.3334

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 4837.5186
$ws.Range("I46").Value = 2163.842
$ws.Range("K46").Value = 2163.842
$ws.Range("M46").Value = -1975.842

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 3815.2
$ws.Range("I100").Value = 3300.6667
$ws.Range("J100").Value = 4035.7144
$ws.Range("K100").Value = 3300.6667
$ws.Range("L100").Value = 4035.7144
$ws.Range("M100").Value = -2759.6667
$ws.Range("N100").Value = -5117.7144

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H133").Value = 68000
$ws.Range("I133").Value = 0
$ws.Range("J133").Value = 68000
$ws.Range("K133").Value = 0
$ws.Range("L133").Value = 68000
$ws.Range("N133").Value = -73060
$ws.Range("M133").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 1526.25
$ws.Range("I122").Value = 1528.7273
$ws.Range("K122").Value = 4586.1819
$ws.Range("L122").Value = 1499
$ws.Range("M122").Value = -2136.1819

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 2138.0952
$ws.Range("I132").Value = 1319
$ws.Range("J132").Value = 4759.2
$ws.Range("K132").Value = 3957
$ws.Range("L132").Value = 14277.6
$ws.Range("M132").Value = -1427
$ws.Range("N132").Value = -19337.6
